$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold a textual "1" (the value must stay text, not become the
# number 1) - force a text format before assigning so Excel keeps it as a
# string instead of silently converting the numeric-looking text to a number.
$textCells = @("C6", "F6", "C9", "C10", "C11", "C12", "C13", "B18", "G30", "B34", "E34", "B37", "E37")
foreach ($addr in $textCells) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = "1"
}

# These cells hold numeric values and stay numeric.
$ws.Range("C8").Value = 45427
$ws.Range("C14").Value = 1
